$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the TEST_NAME comment text to a module-scoped variant
$ws.Range("I2").Value = "測試名稱0"

# Fill in previously-blank flag cells (Nullable = Y, index-position numbers)
$ws.Range("G2").Value = "Y"
$ws.Range("N2").Value = "1"
$ws.Range("M6").Value = "1"
$ws.Range("M7").Value = "2"
$ws.Range("N8").Value = "2"

# Update the active selection / scroll position for the sheet view
$ws.Range("G2").Select()
